# Adding LCA data to optimization: new technology rows (solar collector,
# district-heating boilers, green electricity, CHP units) across the DHW,
# HEATING and ELECTRICITY sheets, plus a couple of formula tweaks on the
# COOLING and ELECTRICITY sheets (accounting for distribution losses / a
# revised PV cost assumption).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# DHW sheet: add "solar collector" (T7 / SC) as row 4
# ---------------------------------------------------------------------------
$dhw = $wb.Worksheets.Item("DHW")

$dhw.Range("A2:H2").Copy()
$dhw.Range("A4:H4").PasteSpecial(-4122)

$dhw.Range("A4").Value = "solar collector"
$dhw.Range("B4").Value = "T7"
$dhw.Range("C4").Value = "SC"
$dhw.Range("D4").Value = 0.7
$dhw.Range("E4").Value = 0.277
$dhw.Range("F4").Value = 0.013
$dhw.Range("G4").Value = 0
$dhw.Range("H4").Value = "KBOB 2019, costs in USD-2015"

# stray formatted (empty) cell that shows up one column over, as in the
# authored workbook
$dhw.Range("A4").Copy()
$dhw.Range("I4").PasteSpecial(-4122)

$dhw.Range("J15").Select()

# ---------------------------------------------------------------------------
# HEATING sheet: add "solar collector" plus three district-heating boiler
# technologies (bio gas / agricultural bio gas / natural gas)
# ---------------------------------------------------------------------------
$heating = $wb.Worksheets.Item("HEATING")

$heating.Range("A2:H2").Copy()
$heating.Range("A3:H6").PasteSpecial(-4122)

$heating.Range("A3").Value = "solar collector"
$heating.Range("B3").Value = "T7"
$heating.Range("C3").Value = "SC"
$heating.Range("D3").Value = 0.7
$heating.Range("E3").Value = 0.277
$heating.Range("F3").Value = 0.013
$heating.Range("G3").Value = 0
$heating.Range("H3").Value = "KBOB 2019, costs in USD-2015"

$heating.Range("A4").Value = "district heating - bio gas-fired boiler"
$heating.Range("B4").Value = "T23"
$heating.Range("C4").Value = "DH"
$heating.Range("D4").ClearContents()
$heating.Range("E4").Value = 0.8307
$heating.Range("F4").Value = 0.106314
$heating.Range("G4").ClearContents()
$heating.Range("H4").Value = "from CEA, costs in USD-2015"

$heating.Range("A5").Value = "district heating - agricultural bio gas-fired boiler"
$heating.Range("B5").Value = "T24"
$heating.Range("C5").Value = "DH"
$heating.Range("D5").ClearContents()
$heating.Range("E5").Value = 0.176514
$heating.Range("F5").Value = 0.0432276
$heating.Range("G5").ClearContents()
$heating.Range("H5").Value = "from CEA, costs in USD-2015"

$heating.Range("A6").Value = "district heating - natural gas-fired boiler"
$heating.Range("B6").Value = "T25"
$heating.Range("C6").Value = "DH"
$heating.Range("D6").ClearContents()
$heating.Range("E6").Value = 0.172614
$heating.Range("F6").Value = 2.7612
$heating.Range("G6").ClearContents()
$heating.Range("H6").Value = "from CEA, costs in USD-2015"

# stray formatted (empty) cells in column I next to rows 4-6
$heating.Range("A4").Copy()
$heating.Range("I4").PasteSpecial(-4122)
$heating.Range("A5").Copy()
$heating.Range("I5").PasteSpecial(-4122)
$heating.Range("A6").Copy()
$heating.Range("I6").PasteSpecial(-4122)

$heating.Range("H13").Select()

# ---------------------------------------------------------------------------
# COOLING sheet: account for extra 10% distribution losses in the district
# cooling network costs, and refresh the "reference" wording to the
# USD-2015 costing basis
# ---------------------------------------------------------------------------
$cooling = $wb.Worksheets.Item("COOLING")

$cooling.Range("G5").Formula = "=0.2/2.7*1.1"
$cooling.Range("G6").Formula = "=(ELECTRICITY!G4/4)*1.1"

$cooling.Range("H3").Value = "Embodied is neglected, only electricity"
$cooling.Range("H4").Value = "Embodied is neglected, only electricity"

$cooling.Range("G5").Select()

# ---------------------------------------------------------------------------
# ELECTRICITY sheet: refresh the reference wording to the USD-2015 costing
# basis, revise the grid-mix cost assumption, and add Green Electricity plus
# Natural gas / Bio gas / Agricultural bio gas CHP technologies
# ---------------------------------------------------------------------------
$electricity = $wb.Worksheets.Item("ELECTRICITY")

$electricity.Range("H3").Value = "ecoinvent 3.4 - electricity production, photovoltaic, 3kWp flat-roof installation, single-Si RoW electricity, medium voltage, costs in USD-2015"

$electricity.Range("G4").Formula = "=0.22*0.75"
$electricity.Range("H4").Value = "ecoinvent 3.4 - market for electricity, medium voltage, SG, costs in USD-2015"

$electricity.Range("A2:H2").Copy()
$electricity.Range("A5:H8").PasteSpecial(-4122)

$electricity.Range("A5").Value = "Green Electricity"
$electricity.Range("B5").Value = "T7"
$electricity.Range("C5").Value = "mix"
$electricity.Range("D5").ClearContents()
$electricity.Range("E5").Value = 0.034
$electricity.Range("F5").Value = 0.004
$electricity.Range("G5").ClearContents()
$electricity.Range("H5").Value = "from CEA, costs in USD-2015"

$electricity.Range("A6").Value = "Natural gas CHP"
$electricity.Range("B6").Value = "T8"
$electricity.Range("C6").Value = "NG"
$electricity.Range("D6").ClearContents()
$electricity.Range("E6").Value = 2.2932
$electricity.Range("F6").Value = 0.14508
$electricity.Range("G6").ClearContents()
$electricity.Range("H6").Value = "from CEA, costs in USD-2015"

$electricity.Range("A7").Value = "Bio gas CHP"
$electricity.Range("B7").Value = "T9"
$electricity.Range("C7").Value = "BG"
$electricity.Range("D7").ClearContents()
$electricity.Range("E7").Value = 0.66378
$electricity.Range("F7").Value = 0.08892
$electricity.Range("G7").ClearContents()
$electricity.Range("H7").Value = "from CEA, costs in USD-2015"

$electricity.Range("A8").Value = "Agricultural Bio gas CHP"
$electricity.Range("B8").Value = "T10"
$electricity.Range("C8").Value = "BG"
$electricity.Range("D8").ClearContents()
$electricity.Range("E8").Value = 0.12168
$electricity.Range("F8").Value = 0.03861
$electricity.Range("G8").ClearContents()
$electricity.Range("H8").Value = "from CEA, costs in USD-2015"

$electricity.Range("H18").Select()

# ---------------------------------------------------------------------------
# leave HEATING as the active sheet/tab, matching the authored workbook
# ---------------------------------------------------------------------------
$heating.Activate()
$heating.Range("H13").Select()
